$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 207.25807
$ws.Range("I8").Value = 62
$ws.Range("K8").Value = 186
$ws.Range("M8").Value = -47

$ws.Range("H132").Value = 20411856
$ws.Range("I132").Value = 27781128
$ws.Range("K132").Value = 83343384
$ws.Range("M132").Value = -83340854

$ws.Range("H138").Value = 7569.7856
$ws.Range("J138").Value = 10622.75
$ws.Range("L138").Value = 31868.25
$ws.Range("N138").Value = -42148.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2242.5544
$ws.Range("I32").Value = 2124.2207
$ws.Range("J32").Value = 2850
$ws.Range("K32").Value = 2124.2207
$ws.Range("L32").Value = 2850
$ws.Range("M32").Value = -1837.2207
$ws.Range("N32").Value = -3424

$ws.Range("H74").Value = 1513
$ws.Range("I74").Value = 1375.25
$ws.Range("K74").Value = 1375.25
$ws.Range("M74").Value = -501.25

$ws.Range("H77").Value = 1513
$ws.Range("I77").Value = 1375.25
$ws.Range("K77").Value = 6876.25
$ws.Range("M77").Value = -2508.25

$ws.Range("H103").Value = 49497.5
$ws.Range("J103").Value = 49497.5
$ws.Range("L103").Value = 49497.5
$ws.Range("N103").Value = -51841.5

$ws.Range("H122").Value = 3856.7844
$ws.Range("I122").Value = 3269.2
$ws.Range("J122").Value = 5993.4546
$ws.Range("K122").Value = 9807.599999999999
$ws.Range("L122").Value = 17980.3638
$ws.Range("M122").Value = -7357.599999999999
$ws.Range("N122").Value = -22880.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 27500
$ws.Range("J56").Value = 40000
$ws.Range("L56").Value = 40000
$ws.Range("N56").Value = -41478

$ws.Range("H107").Value = 4701.75
$ws.Range("I107").Value = 3723.2727
$ws.Range("J107").Value = 8289.5
$ws.Range("K107").Value = 3723.2727
$ws.Range("L107").Value = 8289.5
$ws.Range("M107").Value = -1803.2727
$ws.Range("N107").Value = -12129.5

$ws.Range("H122").Value = 87000
$ws.Range("J122").Value = 87000
$ws.Range("L122").Value = 87000
$ws.Range("N122").Value = -96800

$ws.Range("H134").Value = 1650.2858
$ws.Range("I134").Value = 1416.907
$ws.Range("K134").Value = 4250.721
$ws.Range("M134").Value = -1715.721

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2387.4167
$ws.Range("I16").Value = 2383.6
$ws.Range("K16").Value = 2383.6
$ws.Range("M16").Value = -2096.6

$ws.Range("H31").Value = 1677.0264
$ws.Range("I31").Value = 1500.909
$ws.Range("K31").Value = 1500.909
$ws.Range("M31").Value = -1205.909

$ws.Range("H34").Value = 1677.0264
$ws.Range("I34").Value = 1500.909
$ws.Range("K34").Value = 1500.909
$ws.Range("M34").Value = -1298.909

$ws.Range("H74").Value = 66825.336
$ws.Range("J74").Value = 66825.336
$ws.Range("L74").Value = 66825.336
$ws.Range("N74").Value = -68573.336

$ws.Range("H77").Value = 66825.336
$ws.Range("J77").Value = 66825.336
$ws.Range("L77").Value = 200476.008
$ws.Range("N77").Value = -209212.008

$ws.Range("H112").Value = 99666.664
$ws.Range("I112").Value = 60000
$ws.Range("J112").Value = 119500
$ws.Range("K112").Value = 60000
$ws.Range("L112").Value = 119500
$ws.Range("M112").Value = -58523
$ws.Range("N112").Value = -122454

$ws.Range("H113").Value = 2387.4167
$ws.Range("I113").Value = 2383.6
$ws.Range("K113").Value = 2383.6
$ws.Range("M113").Value = -213.5999999999999

$ws.Range("H132").Value = 1509.2084
$ws.Range("I132").Value = 1305.5
$ws.Range("K132").Value = 3916.5
$ws.Range("M132").Value = -1386.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1057.2727
$ws.Range("I5").Value = 502.33334
$ws.Range("J5").Value = 1723.2
$ws.Range("K5").Value = 1507.00002
$ws.Range("L5").Value = 5169.6
$ws.Range("M5").Value = -1395.00002
$ws.Range("N5").Value = -5393.6

$ws.Range("H13").Value = 463.5
$ws.Range("J13").Value = 435
$ws.Range("L13").Value = 1305
$ws.Range("N13").Value = -1641

$ws.Range("H17").Value = 1182.6364
$ws.Range("I17").Value = 825
$ws.Range("K17").Value = 2475
$ws.Range("M17").Value = -2306

$ws.Range("H129").Value = 3522.7693
$ws.Range("I129").Value = 809.7143
$ws.Range("J129").Value = 6688
$ws.Range("K129").Value = 2429.1429
$ws.Range("L129").Value = 20064
$ws.Range("M129").Value = 2570.8571
$ws.Range("N129").Value = -30064

$ws.Range("N134").ClearContents()
$ws.Range("H134").Value = 4447.25
$ws.Range("I134").Value = 4447.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13341.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8271.75

$ws.Range("H135").Value = 1057.2727
$ws.Range("I135").Value = 502.33334
$ws.Range("J135").Value = 1723.2
$ws.Range("K135").Value = 4521.00006
$ws.Range("L135").Value = 15508.8
$ws.Range("M135").Value = -1986.00006
$ws.Range("N135").Value = -20578.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M20").ClearContents()
$ws.Range("H20").Value = 30233
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0

$ws.Range("N24").ClearContents()
$ws.Range("H24").Value = 100524000
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0

$ws.Range("H122").Value = 4288.684
$ws.Range("I122").Value = 2750.0715
$ws.Range("J122").Value = 8596.799999999999
$ws.Range("K122").Value = 8250.2145
$ws.Range("L122").Value = 25790.4
$ws.Range("M122").Value = -5800.2145
$ws.Range("N122").Value = -30690.4

$ws.Range("H132").Value = 5380.328
$ws.Range("I132").Value = 5203.727
$ws.Range("K132").Value = 15611.181
$ws.Range("M132").Value = -13081.181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5009
$ws.Range("I4").Value = 5009
$ws.Range("K4").Value = 5009
$ws.Range("M4").Value = -4896

$ws.Range("H5").Value = 11800
$ws.Range("I5").Value = 800
$ws.Range("J5").Value = 15466.667
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 15466.667
$ws.Range("M5").Value = -687
$ws.Range("N5").Value = -15692.667

$ws.Range("N26").ClearContents()
$ws.Range("H26").Value = 11210
$ws.Range("I26").Value = 11210
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 11210
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -10915

$ws.Range("H28").Value = 5009
$ws.Range("I28").Value = 5009
$ws.Range("K28").Value = 5009
$ws.Range("M28").Value = -4777

$ws.Range("H37").Value = 5009
$ws.Range("I37").Value = 5009
$ws.Range("K37").Value = 5009
$ws.Range("M37").Value = -4902

$ws.Range("H41").Value = 24111
$ws.Range("I41").Value = 24111
$ws.Range("K41").Value = 24111
$ws.Range("M41").Value = -23673

$ws.Range("H97").Value = 46868.8
$ws.Range("J97").Value = 46868.8
$ws.Range("L97").Value = 46868.8
$ws.Range("N97").Value = -48850.8

$ws.Range("H117").Value = 90001
$ws.Range("J117").Value = 90001
$ws.Range("L117").Value = 90001
$ws.Range("N117").Value = -99179

$ws.Range("H132").Value = 5095.206
$ws.Range("I132").Value = 3376.0715
$ws.Range("K132").Value = 10128.2145
$ws.Range("M132").Value = -7598.2145

$ws.Range("H136").Value = 5037.353
$ws.Range("J136").Value = 4082.8333
$ws.Range("L136").Value = 12248.4999
$ws.Range("N136").Value = -17348.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 10015005
$ws.Range("I28").Value = 15007508
$ws.Range("J28").Value = 29999
$ws.Range("K28").Value = 15007508
$ws.Range("L28").Value = 29999
$ws.Range("M28").Value = -15007160
$ws.Range("N28").Value = -30695

$ws.Range("H40").Value = 33333
$ws.Range("J40").Value = 33333
$ws.Range("L40").Value = 33333
$ws.Range("N40").Value = -33631

$ws.Range("H43").Value = 30014.5
$ws.Range("I43").Value = 29999
$ws.Range("J43").Value = 30030
$ws.Range("K43").Value = 29999
$ws.Range("L43").Value = 30030
$ws.Range("M43").Value = -29850
$ws.Range("N43").Value = -30328

$ws.Range("H112").Value = 79693.5
$ws.Range("J112").Value = 79693.5
$ws.Range("L112").Value = 79693.5
$ws.Range("N112").Value = -82647.5

$ws.Range("N118").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0

$ws.Range("N120").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0

$ws.Range("H126").Value = 10380.091
$ws.Range("I126").Value = 8968.1
$ws.Range("K126").Value = 26904.3
$ws.Range("M126").Value = -24434.3

$ws.Range("H136").Value = 2192.7576
$ws.Range("I136").Value = 1927.8148
$ws.Range("J136").Value = 3385
$ws.Range("K136").Value = 5783.4444
$ws.Range("L136").Value = 10155
$ws.Range("M136").Value = -3233.4444
$ws.Range("N136").Value = -15255
